# This script applies a weekly update to the "Coco" price sheet.
# Rows 5-26 shift down by one record (each row now shows the figures
# that used to belong to the following row), row 5 receives a new
# week of data, and the data that used to sit in row 26 is preserved
# as the new row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 44414; $ws.Range("N5").Value = 25000; $ws.Range("O5").Value = 25000; $ws.Range("P5").Value = 25000; $ws.Range("S5").Value = 1250
$ws.Range("D6").Value = 44356; $ws.Range("M6").Value = 15
$ws.Range("D7").Value = 44396; $ws.Range("M7").Value = 12; $ws.Range("N7").Value = 24000; $ws.Range("O7").Value = 24000; $ws.Range("P7").Value = 24000; $ws.Range("S7").Value = 1200
$ws.Range("D8").Value = 44221; $ws.Range("M8").Value = 30; $ws.Range("N8").Value = 25000; $ws.Range("O8").Value = 25000; $ws.Range("P8").Value = 25000; $ws.Range("S8").Value = 1250
$ws.Range("D9").Value = 44175; $ws.Range("M9").Value = 25; $ws.Range("N9").Value = 23000; $ws.Range("O9").Value = 23000; $ws.Range("P9").Value = 23000; $ws.Range("S9").Value = 1150
$ws.Range("D10").Value = 44363
$ws.Range("D11").Value = 44349; $ws.Range("M11").Value = 30; $ws.Range("N11").Value = 24000; $ws.Range("O11").Value = 24000; $ws.Range("P11").Value = 24000; $ws.Range("S11").Value = 1200
$ws.Range("D12").Value = 44222; $ws.Range("N12").Value = 25000; $ws.Range("O12").Value = 25000; $ws.Range("P12").Value = 25000; $ws.Range("S12").Value = 1250
$ws.Range("D13").Value = 44377; $ws.Range("M13").Value = 15; $ws.Range("N13").Value = 20000; $ws.Range("O13").Value = 20000; $ws.Range("P13").Value = 20000; $ws.Range("S13").Value = 1000
$ws.Range("D14").Value = 44400; $ws.Range("M14").Value = 5; $ws.Range("N14").Value = 24000; $ws.Range("O14").Value = 24000; $ws.Range("P14").Value = 24000; $ws.Range("S14").Value = 1200
$ws.Range("D15").Value = 44194; $ws.Range("M15").Value = 20; $ws.Range("N15").Value = 20000; $ws.Range("O15").Value = 20000; $ws.Range("P15").Value = 20000; $ws.Range("S15").Value = 1000
$ws.Range("D16").Value = 44390; $ws.Range("M16").Value = 10; $ws.Range("N16").Value = 24000; $ws.Range("O16").Value = 24000; $ws.Range("P16").Value = 24000; $ws.Range("S16").Value = 1200
$ws.Range("D17").Value = 44412; $ws.Range("M17").Value = 20
$ws.Range("D18").Value = 44214; $ws.Range("M18").Value = 15
$ws.Range("D19").Value = 44238; $ws.Range("M19").Value = 30
$ws.Range("D20").Value = 44231; $ws.Range("M20").Value = 15; $ws.Range("N20").Value = 25000; $ws.Range("O20").Value = 25000; $ws.Range("P20").Value = 25000; $ws.Range("S20").Value = 1250
$ws.Range("D21").Value = 44391; $ws.Range("M21").Value = 10
$ws.Range("D22").Value = 44389; $ws.Range("M22").Value = 20; $ws.Range("N22").Value = 24000; $ws.Range("O22").Value = 24000; $ws.Range("P22").Value = 24000; $ws.Range("S22").Value = 1200
$ws.Range("D23").Value = 44251; $ws.Range("N23").Value = 25000; $ws.Range("O23").Value = 25000; $ws.Range("P23").Value = 25000; $ws.Range("S23").Value = 1250
$ws.Range("D24").Value = 44382; $ws.Range("N24").Value = 20000; $ws.Range("O24").Value = 20000; $ws.Range("P24").Value = 20000; $ws.Range("S24").Value = 1000
$ws.Range("D25").Value = 44232
$ws.Range("D26").Value = 44398; $ws.Range("M26").Value = 15

# New row 27 (previously row 26 data, before the weekly shift)
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 44334
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108007
$ws.Range("J27").Value = "Coco"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 25000
$ws.Range("Q27").Value = "`$/malla 20 unidades"
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 1250
$ws.Range("T27").Value = 20

# Preserve the date formatting used by the rest of column D
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update the used range dimension to reflect the newly added row
Write-Host "Applied weekly shift; new dimension should be A1:T27"
